$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "optimise the card bag in game shop"
# The old "card bag" shop rows (kabaolanse/kabaohuangse/kabaohongse,
# kapaibugeibao(*), sucaidai*) are replaced by three unified kabao1 items,
# and every subsequent shop-entry row shifts up to close the gap.
# ---------------------------------------------------------------------------

# Final target content for A4:C35 (row, A, B-text, C).
# Rows 1-3 are header rows and their visible text does not change in this
# edit (only the internal shared-string index shifts elsewhere in the
# table, which Excel recomputes on save) - so they are left untouched.
$data = @(
    @(4, 15000001, "kabao1",              1),
    @(5, 15000002, "kabao1v2",            1),
    @(6, 15000003, "kabao1v3",            1),
    @(7, 15000014, "suijihuanshouka",     2),
    @(8, 15000015, "suijiwuqika",         2),
    @(9, 15000016, "suijimofaka",         2),
    @(10,15000017, "fuwen-aier",          2),
    @(11,15000018, "fuwen-puer",          2),
    @(12,15000019, "fuwen-chamu",         2),
    @(13,15000020, "zhongxinghuoliyaoji", 2),
    @(14,15000021, "daxinghuoliyaoji",    2),
    @(15,15000022, "zhongxingmofayaoji",  2),
    @(16,15000023, "daxingmofayaoji",     2),
    @(17,15000024, "zhongxingtiliyaoji",  2),
    @(18,15000025, "daxingtiliyaoji",     2),
    @(19,15000026, "gangtiexiulichui",    2),
    @(20,15000027, "shenshengxiulichui",  2),
    @(21,15000028, "jingyanzhishu",       3),
    @(22,15000029, "nenliangzhishu",      3),
    @(23,15000030, "binggan",             3),
    @(24,15000031, "hongsejiaonan",       3),
    @(25,15000032, "lansejiaonan",        3),
    @(26,15000033, "shuijingqiu",         3),
    @(27,15000034, "zuoqiheibao",         3),
    @(28,15000035, "zuoqiying",           3),
    @(29,15000036, "yaoshuistr",          3),
    @(30,15000037, "yaoshuiintl",         3),
    @(31,15000038, "yaoshui1",            3),
    @(32,15000039, "yaoshui2",            3),
    @(33,15000040, "yaoshuiagi",          3),
    @(34,15000041, "yaoshuiperc",         3),
    @(35,15000042, "yaoshuiendu",         3)
)

# Rows (1-based, matching the list above) whose B cell must carry the
# bordered "kabao card" style (same formatting already used by the last
# rows of the table, e.g. B39 in the original sheet).
$styledRows = @(29,30,31,32,33,34,35)

# Grab a style source BEFORE we touch anything - row 39 col B already has
# the target style (s="3") in the original workbook.
$styleSource = $ws.Range("B39")

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
}

# Apply the bordered style to the newly relocated "card bag" rows.
$styleSource.Copy()
foreach ($r in $styledRows) {
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
}

# Drop the now-unused tail rows (36-45) - this also shrinks the table /
# autofilter / dimension to A1:C35 automatically.
$ws.Rows("36:45").Delete()

# Reset the scroll position back to the top-left (no more topLeftCell="A19").
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
